$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear out column C entirely (it is being removed) ---
$ws.Range("C1:C4").Clear()

# --- Row 2 now holds the "2*theta/°" / "R(30kV)/Imp/s" headers ---
$ws.Range("A2").Value = "2*theta/°"
$ws.Range("B2").Value = "R(30kV)/Imp/s"

# --- Clear the old single data row (row 4) before re-populating ---
$ws.Range("A4:B4").Clear()

# --- New data rows 4-9 ---
$ws.Range("A4").Value = 10
$ws.Range("B4").Value = 13347

$ws.Range("A5").Value = 30
$ws.Range("B5").Value = 13299

$ws.Range("A6").Value = 50
$ws.Range("B6").Value = 12012

$ws.Range("A7").Value = 70
$ws.Range("B7").Value = 1073

$ws.Range("A8").Value = 90
$ws.Range("B8").Value = 3535

$ws.Range("A9").Value = 110
$ws.Range("B9").Value = 1084

# --- Update the visible selection to match the new used range ---
$ws.Range("A1:B9").Select() | Out-Null
